# Add a new "intervention_type" column (K) to the clinical trials list sheet,
# populating the header and values for each row, matching the style of the
# existing header row for the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell
$ws.Range("K1").Value = "intervention_type"
# Match the formatting (bold, centered, bordered) used by the other header cells
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for K2:K17
$ws.Range("K2").Value  = "OTHER"
$ws.Range("K3").Value  = "OTHER"
$ws.Range("K4").Value  = "OTHER"
$ws.Range("K5").Value  = "OTHER"
$ws.Range("K6").Value  = "BEHAVIORAL"
$ws.Range("K7").Value  = "OTHER"
$ws.Range("K8").Value  = "DIETARY_SUPPLEMENT"
$ws.Range("K9").Value  = "OTHER"
$ws.Range("K10").Value = "OTHER"
$ws.Range("K11").Value = "OTHER"

# Row 12 has no known intervention type; leave the cell present but blank,
# mirroring the other blank cells already found on that row (e.g. C12).
$ws.Range("C12").Copy($ws.Range("K12"))

$ws.Range("K13").Value = "DIETARY_SUPPLEMENT"
$ws.Range("K14").Value = "DIETARY_SUPPLEMENT"
$ws.Range("K15").Value = "DIETARY_SUPPLEMENT"
$ws.Range("K16").Value = "DIETARY_SUPPLEMENT"
$ws.Range("K17").Value = "OTHER"
